$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 866.6667
$ws.Range("I80").Value = 833.3333
$ws.Range("J80").Value = 900
$ws.Range("K80").Value = 2499.9999
$ws.Range("L80").Value = 2700
$ws.Range("M80").Value = -1501.9999
$ws.Range("N80").Value = -4696

$ws.Range("H83").Value = 866.6667
$ws.Range("I83").Value = 833.3333
$ws.Range("J83").Value = 900
$ws.Range("K83").Value = 7499.9997
$ws.Range("L83").Value = 8100
$ws.Range("M83").Value = -2507.9997
$ws.Range("N83").Value = -18084

$ws.Range("H96").Value = 546.5
$ws.Range("I96").Value = 429.5
$ws.Range("J96").Value = 1014.5
$ws.Range("K96").Value = 1288.5
$ws.Range("L96").Value = 3043.5
$ws.Range("M96").Value = 84.5
$ws.Range("N96").Value = -5789.5

$ws.Range("H97").Value = 7211.3335
$ws.Range("J97").Value = 7211.3335
$ws.Range("L97").Value = 21634.0005
$ws.Range("N97").Value = -22626.0005

$ws.Range("H137").Value = 3707.5
$ws.Range("I137").Value = 3280.75
$ws.Range("J137").Value = 3897.1667
$ws.Range("K137").Value = 9842.25
$ws.Range("L137").Value = 11691.5001
$ws.Range("M137").Value = -7292.25
$ws.Range("N137").Value = -16791.5001

$ws.Range("H141").Value = 770.5
$ws.Range("I141").Value = 770.5
$ws.Range("K141").Value = 2311.5
$ws.Range("M141").Value = 2868.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2440.4
$ws.Range("I2").Value = 2550.5
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 2550.5
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -2437.5
$ws.Range("N2").Value = -2226

$ws.Range("H32").Value = 13369.9
$ws.Range("I32").Value = 8954.75
$ws.Range("J32").Value = 19992.625
$ws.Range("K32").Value = 8954.75
$ws.Range("L32").Value = 19992.625
$ws.Range("M32").Value = -8667.75
$ws.Range("N32").Value = -20566.625

$ws.Range("H74").Value = 1990
$ws.Range("I74").Value = 1990
$ws.Range("K74").Value = 1990
$ws.Range("M74").Value = -1116

$ws.Range("H77").Value = 1990
$ws.Range("I77").Value = 1990
$ws.Range("K77").Value = 9950
$ws.Range("M77").Value = -5582

$ws.Range("H116").Value = 2440.4
$ws.Range("I116").Value = 2550.5
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 2550.5
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = -256.5
$ws.Range("N116").Value = -6588

$ws.Range("H132").Value = 5082.2666
$ws.Range("I132").Value = 5082.2666
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 15246.7998
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -12716.7998
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2440.4
$ws.Range("I3").Value = 2550.5
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 2550.5
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -2436.5
$ws.Range("N3").Value = -2228

$ws.Range("H22").Value = 221.8
$ws.Range("I22").Value = 117.85714
$ws.Range("K22").Value = 117.85714
$ws.Range("M22").Value = 55.14286

$ws.Range("H80").Value = 1362.7
$ws.Range("I80").Value = 903.25
$ws.Range("K80").Value = 903.25
$ws.Range("M80").Value = 94.75

$ws.Range("H83").Value = 1362.7
$ws.Range("I83").Value = 903.25
$ws.Range("K83").Value = 4516.25
$ws.Range("M83").Value = 475.75

$ws.Range("H86").Value = 4582.846
$ws.Range("I86").Value = 2646.3333
$ws.Range("K86").Value = 2646.3333
$ws.Range("M86").Value = -1523.3333

$ws.Range("H89").Value = 4582.846
$ws.Range("I89").Value = 2646.3333
$ws.Range("K89").Value = 13231.6665
$ws.Range("M89").Value = -7615.666499999999

$ws.Range("H134").Value = 3001
$ws.Range("I134").Value = 3001
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 9003
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -6468
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1112.5
$ws.Range("I16").Value = 1150
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 1150
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -863
$ws.Range("N16").Value = -1574

$ws.Range("H25").Value = 15000
$ws.Range("I25").Value = 15000
$ws.Range("K25").Value = 15000
$ws.Range("M25").Value = -14826

$ws.Range("H31").Value = 20276.215
$ws.Range("I31").Value = 11297.546
$ws.Range("K31").Value = 11297.546
$ws.Range("M31").Value = -11002.546

$ws.Range("H34").Value = 20276.215
$ws.Range("I34").Value = 11297.546
$ws.Range("K34").Value = 11297.546
$ws.Range("M34").Value = -11095.546

$ws.Range("H58").Value = 4293.75
$ws.Range("I58").Value = 4293.75
$ws.Range("K58").Value = 4293.75
$ws.Range("M58").Value = -4090.75

$ws.Range("H113").Value = 1112.5
$ws.Range("I113").Value = 1150
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1150
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1020
$ws.Range("N113").Value = -5340

$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -54920

$ws.Range("H134").Value = 1011
$ws.Range("I134").Value = 870.1667
$ws.Range("J134").Value = 1292.6666
$ws.Range("K134").Value = 2610.5001
$ws.Range("L134").Value = 3877.9998
$ws.Range("M134").Value = -75.5001000000002
$ws.Range("N134").Value = -8947.9998

$ws.Range("H136").Value = 4293.75
$ws.Range("I136").Value = 4293.75
$ws.Range("K136").Value = 12881.25
$ws.Range("M136").Value = -10331.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 279
$ws.Range("I97").Value = 450
$ws.Range("K97").Value = 1350
$ws.Range("M97").Value = -854

$ws.Range("H113").Value = 1658
$ws.Range("I113").Value = 975
$ws.Range("K113").Value = 2925
$ws.Range("M113").Value = -755

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H102").Value = 3742.2856
$ws.Range("I102").Value = 3742.2856
$ws.Range("K102").Value = 3742.2856
$ws.Range("M102").Value = -2120.2856

$ws.Range("H122").Value = 8093
$ws.Range("I122").Value = 2994.875
$ws.Range("K122").Value = 8984.625
$ws.Range("M122").Value = -6534.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 618.36365
$ws.Range("I16").Value = 607.7778
$ws.Range("K16").Value = 607.7778
$ws.Range("M16").Value = -437.7778

$ws.Range("H20").Value = 53333.332
$ws.Range("J20").Value = 53333.332
$ws.Range("L20").Value = 53333.332
$ws.Range("N20").Value = -53785.332

$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

$ws.Range("H40").Value = 7500
$ws.Range("I40").Value = 7000
$ws.Range("K40").Value = 7000
$ws.Range("M40").Value = -6864

$ws.Range("H46").Value = 919.8
$ws.Range("I46").Value = 800
$ws.Range("K46").Value = 800
$ws.Range("M46").Value = -612

$ws.Range("H122").Value = 6299
$ws.Range("I122").Value = 5973
$ws.Range("J122").Value = 6625
$ws.Range("K122").Value = 17919
$ws.Range("L122").Value = 19875
$ws.Range("M122").Value = -15469
$ws.Range("N122").Value = -24775

$ws.Range("H132").Value = 29184.125
$ws.Range("I132").Value = 26210.428
$ws.Range("K132").Value = 78631.284
$ws.Range("M132").Value = -76101.284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3658.3333
$ws.Range("I62").Value = 3487.5
$ws.Range("K62").Value = 3487.5
$ws.Range("M62").Value = -2863.5

$ws.Range("H65").Value = 3658.3333
$ws.Range("I65").Value = 3487.5
$ws.Range("K65").Value = 17437.5
$ws.Range("M65").Value = -14317.5

$ws.Range("H81").Value = 2000
$ws.Range("J81").Value = 2000
$ws.Range("L81").Value = 4000
$ws.Range("N81").Value = -6122

$ws.Range("H84").Value = 2000
$ws.Range("J84").Value = 2000
$ws.Range("L84").Value = 20000
$ws.Range("N84").Value = -30608

$ws.Range("H113").Value = 490.33334
$ws.Range("I113").Value = 474.4
$ws.Range("J113").Value = 570
$ws.Range("K113").Value = 1423.2
$ws.Range("L113").Value = 1710
$ws.Range("M113").Value = 746.8000000000002
$ws.Range("N113").Value = -6050

$ws.Range("H132").Value = 4558.4443
$ws.Range("I132").Value = 3603.4666
$ws.Range("J132").Value = 9333.333000000001
$ws.Range("K132").Value = 10810.3998
$ws.Range("L132").Value = 27999.999
$ws.Range("M132").Value = -8280.399800000001
$ws.Range("N132").Value = -33059.999

$ws.Range("H136").Value = 2451.625
$ws.Range("I136").Value = 2166.2144
$ws.Range("J136").Value = 4449.5
$ws.Range("K136").Value = 6498.6432
$ws.Range("M136").Value = -3948.6432
$ws.Range("N136").Value = -18448.5
